# Optuna Attempt (go back with original)
# Update forecast values on the "Forecast Comparison" sheet and the
# corresponding roll-up totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------------

# Row 2 (W8)
$wsForecast.Range("D2").Value = 7
$wsForecast.Range("H2").Value = 14.1
$wsForecast.Range("L2").Value = 0.82

# Row 3 (W9)
$wsForecast.Range("H3").Value = 11.23
$wsForecast.Range("L3").Value = 1.1

# Row 4 (W10)
$wsForecast.Range("H4").Value = 8.949999999999999
$wsForecast.Range("L4").Value = 0.96

# Row 5 (W11)
$wsForecast.Range("H5").Value = 7.95
$wsForecast.Range("L5").Value = 0.96

# Row 6 (W12)
$wsForecast.Range("H6").Value = 6.95
$wsForecast.Range("L6").Value = 1.2

# Row 7 (W13)
$wsForecast.Range("H7").Value = 6.8
$wsForecast.Range("L7").Value = 1.16

# Row 8 (W14)
$wsForecast.Range("H8").Value = 5.08
$wsForecast.Range("L8").Value = 0.95

# Row 9 (W15)
$wsForecast.Range("H9").Value = 4.08
$wsForecast.Range("L9").Value = 1.18

# Row 10 (W16)
$wsForecast.Range("H10").Value = 3.52
$wsForecast.Range("L10").Value = 0.98

# Row 11 (W17)
$wsForecast.Range("H11").Value = 2.52
$wsForecast.Range("L11").Value = 1.02

# Row 12 (W18)
$wsForecast.Range("D12").Value = 10
$wsForecast.Range("H12").Value = 1.33
$wsForecast.Range("L12").Value = 0.97

# Row 13 (W19)
$wsForecast.Range("H13").Value = 0.33
$wsForecast.Range("I13").Value = "High"
$wsForecast.Range("L13").Value = 1.02

# Row 14 (W20)
$wsForecast.Range("L14").Value = 0.87

# Row 15 (W21)
$wsForecast.Range("L15").Value = 0.85

# Row 16 (W22)
$wsForecast.Range("L16").Value = 1.07

# Row 17 (W23)
$wsForecast.Range("L17").Value = 1.17

# --- Summary sheet --------------------------------------------------------
# These "numeric-looking" values are stored as text in the workbook, so we
# force a text number format before assigning them to avoid Excel silently
# converting them to real numbers.

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "153"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "78"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "38"

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "8"
